$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.178962333333333
$ws.Range("H2").Value = 3.536887
$ws.Range("I2").Value = 0.001182125215344215
$ws.Range("J2").Value = 0.001182125215344214
$ws.Range("M2").Value = 0.00535
$ws.Range("N2").Value = 0.01605
$ws.Range("O2").Value = 0.003591913026022235
$ws.Range("P2").Value = 0.003591913026022235
$ws.Range("Q2").Value = 0.006307448483333334
$ws.Range("R2").Value = 0.05676703634999999
$ws.Range("S2").Value = 0.000004246090959384224
$ws.Range("T2").Value = 0.000004246090959384223
$ws.Range("G3").Value = 1.178962333333333
$ws.Range("H3").Value = 3.536887
$ws.Range("I3").Value = 0.001182125215344215
$ws.Range("J3").Value = 0.001182125215344214
$ws.Range("O3").Value = 0.9964080869739778
$ws.Range("P3").Value = 0.9964080869739778
$ws.Range("Q3").Value = 1.749706251636333
$ws.Range("R3").Value = 15.747356264727
$ws.Range("S3").Value = 0.001177879124384831
$ws.Range("T3").Value = 0.00117787912438483
$ws.Range("I4").Value = 0.9532080272144655
$ws.Range("J4").Value = 0.9532080272144653
$ws.Range("M4").Value = 0.00535
$ws.Range("N4").Value = 0.01605
$ws.Range("O4").Value = 0.003591913026022235
$ws.Range("P4").Value = 0.003591913026022235
$ws.Range("Q4").Value = 5.086018340116667
$ws.Range("R4").Value = 45.77416506105
$ws.Range("S4").Value = 0.003423840329460596
$ws.Range("T4").Value = 0.003423840329460595
$ws.Range("I5").Value = 0.9532080272144655
$ws.Range("J5").Value = 0.9532080272144653
$ws.Range("O5").Value = 0.9964080869739778
$ws.Range("P5").Value = 0.9964080869739778
$ws.Range("Q5").Value = 1410.877648728135
$ws.Range("S5").Value = 0.9497841868850049
$ws.Range("T5").Value = 0.9497841868850047
$ws.Range("G6").Value = 45.48781433333333
$ws.Range("I6").Value = 0.04560984757019037
$ws.Range("J6").Value = 0.04560984757019036
$ws.Range("M6").Value = 0.00535
$ws.Range("N6").Value = 0.01605
$ws.Range("O6").Value = 0.003591913026022235
$ws.Range("P6").Value = 0.003591913026022235
$ws.Range("Q6").Value = 0.2433598066833333
$ws.Range("R6").Value = 2.19023826015
$ws.Range("S6").Value = 0.0001638266056022554
$ws.Range("T6").Value = 0.0001638266056022554
$ws.Range("G7").Value = 45.48781433333333
$ws.Range("I7").Value = 0.04560984757019037
$ws.Range("J7").Value = 0.04560984757019036
$ws.Range("O7").Value = 0.9964080869739778
$ws.Range("P7").Value = 0.9964080869739778
$ws.Range("Q7").Value = 67.50878366680031
$ws.Range("R7").Value = 607.5790530012029
$ws.Range("S7").Value = 0.04544602096458811
$ws.Range("T7").Value = 0.04544602096458811
